$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 422; existing rows 422..482 shift down to 423..483.
$ws.Rows.Item(422).EntireRow.Insert()

# Populate the newly inserted row 422 with the new data record.
$ws.Range("A422").Value = 3
$ws.Range("B422").Value = "Femacal de La Calera"
$ws.Range("C422").Value = "Coquimbo"
$ws.Range("D422").Value = 44984
$ws.Range("E422").Value = 5
$ws.Range("F422").Value = 100112012
$ws.Range("G422").Value = "Espinaca"
$ws.Range("H422").Value = "Sin especificar"
$ws.Range("I422").Value = "Primera"
$ws.Range("J422").Value = 140
$ws.Range("K422").Value = 6000
$ws.Range("L422").Value = 6500
$ws.Range("M422").Value = 6286
$ws.Range("N422").Value = '$/docena de atados (3 kilos)'
$ws.Range("O422").Value = "Provincia de Quillota"
$ws.Range("P422").Value = 2095
$ws.Range("Q422").Value = 3
$ws.Range("R422").Value = "Hortaliza"
